$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.093
$ws.Range("C4").Value = -11.962
$ws.Range("B7").Value = 5.553000000000001
$ws.Range("A8").Value = -22.191
$ws.Range("A10").Value = -21.053
$ws.Range("D10").Value = -8.068999999999999
$ws.Range("C11").Value = -12.126
$ws.Range("A12").Value = -21.618
$ws.Range("D12").Value = -7.644999999999999
$ws.Range("D13").Value = -8.404
$ws.Range("B14").Value = 5.431
$ws.Range("C14").Value = -12.94
$ws.Range("D14").Value = -7.9
$ws.Range("B15").Value = 5.235
$ws.Range("E16").Value = 17.076
$ws.Range("A18").Value = -21.231
$ws.Range("B18").Value = 8.115
$ws.Range("C18").Value = -11.957
$ws.Range("E18").Value = 17.436
$ws.Range("C19").Value = -11.752
$ws.Range("B20").Value = 7.098999999999999
$ws.Range("C21").Value = -11.881
$ws.Range("E21").Value = 16.778
$ws.Range("E22").Value = 16.503
$ws.Range("A25").Value = -21.898
$ws.Range("E26").Value = 17.128
$ws.Range("C27").Value = -12.162
$ws.Range("E27").Value = 17.141
$ws.Range("B29").Value = 5.205
$ws.Range("D29").Value = -7.666000000000001
$ws.Range("B30").Value = 5.366000000000001
$ws.Range("B31").Value = 5.111
$ws.Range("C31").Value = -12.872
$ws.Range("D32").Value = -8.363
$ws.Range("B35").Value = 8.379
$ws.Range("D35").Value = -7.781000000000001
$ws.Range("A37").Value = -20.249
$ws.Range("C38").Value = -13.052
$ws.Range("E39").Value = 16.817
$ws.Range("B40").Value = 8.382000000000001
$ws.Range("C42").Value = -12.024
$ws.Range("D43").Value = -8.767999999999999
$ws.Range("B44").Value = 4.915999999999999
$ws.Range("C44").Value = -12.865
$ws.Range("E44").Value = 17.129
$ws.Range("C47").Value = -12.029
$ws.Range("D48").Value = -7.689
$ws.Range("D49").Value = -8.372
$ws.Range("B50").Value = 5.211
$ws.Range("D50").Value = -8.085999999999999
$ws.Range("D51").Value = -8.379000000000001
$ws.Range("E51").Value = 16.988
$ws.Range("B54").Value = 4.761
$ws.Range("E54").Value = 16.644
$ws.Range("A55").Value = -21.756
$ws.Range("C56").Value = -12.841
$ws.Range("D56").Value = -8.068
$ws.Range("E57").Value = 16.546
$ws.Range("C58").Value = -12.951
$ws.Range("E58").Value = 16.865
$ws.Range("E60").Value = 17.048
$ws.Range("D61").Value = -7.861999999999999
$ws.Range("E63").Value = 17.61
$ws.Range("C65").Value = -12.567
$ws.Range("A68").Value = -21.526
$ws.Range("B68").Value = 5.657999999999999
$ws.Range("D69").Value = -7.846999999999999
$ws.Range("D71").Value = -7.645
$ws.Range("C73").Value = -12.646
$ws.Range("B76").Value = 5.532
$ws.Range("A77").Value = -19.913
$ws.Range("E77").Value = 16.918
$ws.Range("A78").Value = -20.207
$ws.Range("A79").Value = -20.858
$ws.Range("D79").Value = -7.914
$ws.Range("A80").Value = -20.175
$ws.Range("A81").Value = -21.818
$ws.Range("D81").Value = -7.877000000000001
$ws.Range("A82").Value = -21.969
$ws.Range("E83").Value = 16.884
$ws.Range("A84").Value = -21.925
$ws.Range("E85").Value = 17.051
$ws.Range("E86").Value = 16.626
$ws.Range("B87").Value = 4.853000000000001
$ws.Range("B88").Value = 5.188000000000001
$ws.Range("C90").Value = -12.876
$ws.Range("B92").Value = 4.647
$ws.Range("C92").Value = -12.227
$ws.Range("D92").Value = -8.085000000000001
$ws.Range("C94").Value = -10.902
$ws.Range("C95").Value = -12.001
$ws.Range("B96").Value = 6.290999999999999
$ws.Range("E96").Value = 16.72
$ws.Range("B98").Value = 5.385000000000001
$ws.Range("E98").Value = 16.949
$ws.Range("A101").Value = -21.075
$ws.Range("B101").Value = 5.912999999999999
$ws.Range("C101").Value = -12.45
$ws.Range("A102").Value = -21.358
$ws.Range("B102").Value = 6.813000000000001
